# "Sửa thông tin cá nhân" (Edit personal information)
# Updates the first member's (row 2, Key=1) record and fixes the
# "Đơn vị" (unit) for three other members who belong to "Ban Đào Tạo".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Key = 1, "Ngô Xuân Hinh") — fill in missing fields and correct
# the address / clear the outdated major.
$ws.Range("B2").Value = "N/A"               # LabID
$ws.Range("G2").Value = "N/A"               # Thế hệ (generation)
$ws.Range("H2").Value = ""                  # Chuyên ngành (major) - cleared
$ws.Range("L2").Value = "Hải Dương "        # Địa chỉ (address)
$ws.Range("M2").Value = "N/A"               # Đơn vị (unit)
$ws.Range("N2").Value = "N/A"               # Chức vụ (position)

# Correct "Đơn vị" from "Chưa có" to "Ban Đào Tạo" for rows 18, 21, 24
$ws.Range("M18").Value = "Ban Đào Tạo"
$ws.Range("M21").Value = "Ban Đào Tạo"
$ws.Range("M24").Value = "Ban Đào Tạo"
